$d = $word.ActiveDocument

# 1. Shorten the university name on the cover/title pages:
#    "TRƯỜNG ĐẠI HỌC CẦN THƠ" -> "ĐẠI HỌC CẦN THƠ" (appears twice).
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("TRƯỜNG ĐẠI HỌC CẦN THƠ", $false, $false, $false, $false, $false, `
              $true, 1, $false, "ĐẠI HỌC CẦN THƠ", 2)

# 2. Remove the stray "_GoBack" bookmark left over from the last cursor
#    position (Word drops this automatically on a clean save).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
